# "funding requests now considers min cash rules"
#
# 1. Account_Min_Max sheet gets a new "Category" column inserted right
#    after "Account" (old col B..F shift to C..G), and the two existing
#    min/max rules (A11 min value, A13 min percent) are tagged with the
#    "Cash/MMKT" category so funding requests can respect a minimum cash
#    rule per account.
# 2. The active sheet/selection moves from Fund_Accounts back to Accounts.

$wb = $excel.ActiveWorkbook

$wsMinMax = $wb.Worksheets.Item("Account_Min_Max")

# Insert a new column before column B ("Minimum Value" -> shifts to C, etc.)
$wsMinMax.Columns.Item(2).Insert()

# New column inherits the width the old "Minimum Value" column (now C) had.
$wsMinMax.Columns.Item(2).ColumnWidth = 13.1666666666667

# New header + values for the inserted "Category" column.
$wsMinMax.Range("B1").Value = "Category"
$wsMinMax.Range("B2").Value = "Cash/MMKT"
$wsMinMax.Range("B3").Value = "Cash/MMKT"

# Remember where the cursor was left on this sheet.
$wsMinMax.Range("D18").Select()

# Move focus back to the Accounts sheet (was Fund_Accounts).
$wsAccounts = $wb.Worksheets.Item("Accounts")
$wsAccounts.Activate()
$wsAccounts.Range("F15").Select()
